$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $value + '"'
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue 'D2' '30.194.75'
Set-TextValue 'E2' '  -1.07%  '
Set-TextValue 'D3' '1.836.69'
Set-TextValue 'E3' '  -1.88%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '231.77'
Set-TextValue 'E5' '  -1.67%  '
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'D7' '0.4653'
Set-TextValue 'E7' '  -3.77%  '
Set-TextValue 'D8' '0.2692'
Set-TextValue 'E8' '  -4.10%  '
Set-TextValue 'D9' '0.06265'
Set-TextValue 'E9' '  -3.77%  '
Set-TextValue 'D10' '1.835.37'
Set-TextValue 'E10' '  -1.89%  '
Set-TextValue 'D11' '0.07373'
Set-TextValue 'E11' '  -0.91%  '
Set-TextValue 'D12' '16.06'
Set-TextValue 'E12' '  -1.79%  '
Set-TextValue 'D13' '4.905'
Set-TextValue 'E13' '  -3.28%  '
Set-TextValue 'D14' '83.36'
Set-TextValue 'E14' '  -4.61%  '
Set-TextValue 'D15' '0.6147'
Set-TextValue 'E15' '  -5.12%  '
Set-TextValue 'D16' '30.121.84'
Set-TextValue 'E16' '  -1.16%  '
Set-TextValue 'E17' '  +0.04%  '
Set-TextValue 'D18' '230.31'
Set-TextValue 'E18' '  -0.45%  '
Set-TextValue 'D19' '0.000007256'
Set-TextValue 'E19' '  -4.09%  '
Set-TextValue 'D20' '1.001'
Set-TextValue 'E20' '  +0.09%  '
Set-TextValue 'D21' '12.30'
Set-TextValue 'E21' '  -5.60%  '
Set-TextValue 'D22' '4.842'
Set-TextValue 'E22' '  -6.32%  '
Set-TextValue 'D23' '5.816'
Set-TextValue 'E23' '  -5.12%  '
Set-TextValue 'D24' '9.186'
Set-TextValue 'E24' '  -1.81%  '
Set-TextValue 'D25' '164.97'
Set-TextValue 'E25' '  -1.10%  '
Set-TextValue 'D26' '17.68'
Set-TextValue 'E26' '  -3.80%  '
Set-TextValue 'D27' '1.867'
Set-TextValue 'E27' '  -3.05%  '
Set-TextValue 'D28' '0.1025'
Set-TextValue 'E28' '  -1.36%  '
Set-TextValue 'E29' '  -0.39%  '
Set-TextValue 'D30' '4.048'
Set-TextValue 'E30' '  -5.33%  '
Set-TextValue 'D31' '3.783'
Set-TextValue 'E31' '  -5.43%  '
Set-TextValue 'E32' '  -4.51%  '
Set-TextValue 'D33' '1.134'
Set-TextValue 'E33' '  -4.18%  '
Set-TextValue 'D34' '0.7085'
Set-TextValue 'E34' '  -3.86%  '
Set-TextValue 'D35' '2.718'
Set-TextValue 'E35' '  +0.16%  '
Set-TextValue 'D36' '0.01828'
Set-TextValue 'E36' '  -3.93%  '
Set-TextValue 'D37' '2.643'
Set-TextValue 'E37' '  +0.15%  '
Set-TextValue 'D38' '0.8894'
Set-TextValue 'E38' '  -2.77%  '
Set-TextValue 'D39' '1.932'
Set-TextValue 'E39' '  -5.78%  '
Set-TextValue 'E40' '  +0.36%  '
Set-TextValue 'D41' '103.73'
Set-TextValue 'E41' '  -2.22%  '
Set-TextValue 'D42' '5.453'
Set-TextValue 'E42' '  -2.85%  '
Set-TextValue 'D43' '0.3993'
Set-TextValue 'E43' '  -5.41%  '
Set-TextValue 'D44' '6.936'
Set-TextValue 'E44' '  -4.95%  '
Set-TextValue 'B45' 'Aave'
Set-TextValue 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '59.70'
Set-TextValue 'E45' '  -5.37%  '
Set-TextValue 'B46' 'Algorand'
Set-TextValue 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D46' '0.1186'
Set-TextValue 'E46' '  -4.91%  '
Set-TextValue 'D47' '8.513'
Set-TextValue 'E47' '  -4.77%  '
Set-TextValue 'D48' '0.05514'
Set-TextValue 'E48' '  -2.45%  '
Set-TextValue 'D49' '32.34'
Set-TextValue 'E49' '  -4.05%  '
Set-TextValue 'D50' '1.350'
Set-TextValue 'E50' '  -6.67%  '
Set-TextValue 'D51' '0.3624'
Set-TextValue 'E51' '  -4.80%  '
